$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "icon" column (G) with header and values
$ws.Range("G1").Value = "icon"
$ws.Range("G2").Value = "gg-pin-blue"
$ws.Range("G3").Value = "gg-pin-pink"
$ws.Range("G4").Value = "gg-pin-blue"
$ws.Range("G5").Value = "gg-pin-green"
$ws.Range("G6").Value = "gg-pin-blue"

# Set the column width to match a "best fit" style autosize
$ws.Columns.Item(7).ColumnWidth = 10.5

# Update the selected cell in the sheet view
$ws.Range("H5").Select()
